$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Data")

# The MainNavigationBarTest suite's pages change too often, so a new
# "What" test case (row 47) is inserted to assert on the URL instead.
# This shifts the remaining MainNavigationBarTest / MainNavSubCategoryTest
# rows down by one.
$ws.Rows.Item(47).Insert()

# Shorten "Guest Bartender" to "Guest" for the (now) row 54 entry in the
# MainNavigationBarTest block only.
$ws.Cells.Item(54, 2).Value = "Guest"

# Populate the newly inserted row with the new "What" test case.
$ws.Cells.Item(47, 1).Value = "Y"
$ws.Cells.Item(47, 2).Value = "What"

# Reflect the new active cell/selection.
$ws.Range("B47").Select()
